$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data table (header row + 10 policy rows), columns A..N
$data = @(
    @("", "`$ bold('All')", "`$ bold('Europe')", "France", "Germany", "Italy", "Poland", "Spain", "United Kingdom", "Switzerland", "Japan", "Russia", "Saudi Arabia", "USA"),
    @("Minimum tax of 2% on billionaires'`nwealth, in voluntary countries", 0.605583686875599, 0.675828047921071, 0.723633675105599, 0.697196274981195, 0.697286079574881, 0.592065792846304, 0.663111648704242, 0.679593714571414, 0.617678362283104, 0.468266623747489, 0.573952839504226, 0.624277119259104, 0.595441917895224),
    @("Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates", 0.540988093608868, 0.585623678646934, 0.539457793802665, 0.579485186187685, 0.692316490221985, 0.472797466690091, 0.591471918081834, 0.603676559780572, 0.593361024155999, 0.380169621797873, 0.572838020267568, 0.672003576449909, 0.528398111147994),
    @("L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage", 0.528944847406596, 0.561663143058492, 0.529949341399934, 0.533956166531526, 0.652215636561299, 0.517677955262585, 0.627598223000534, 0.55648566940254, 0.500082132443895, 0.374045379462791, 0.594965697314265, 0.712886788734677, 0.49932861507328),
    @("Debt relief for vulnerable countries, suspending`npayments until they are more able to repay", 0.487204050910086, 0.516560958421424, 0.46984726554592, 0.449136424403242, 0.620418621743099, 0.519362173741282, 0.565471420059344, 0.523374569882964, 0.490341316203245, 0.332814324343931, 0.508667595717108, 0.685180395151638, 0.4655618778736),
    @("Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales", 0.478171616258382, 0.548978153629316, 0.572441643469511, 0.545688255030839, 0.658945367609779, 0.446658827175329, 0.524794616242944, 0.556188134759176, 0.461310668853147, 0.364008514246129, 0.341022488515294, 0.567119980607023, 0.458128546515251),
    @("At least 0.7% of developed countries' GDP in foreign aid", 0.477487340905981, 0.501409443269908, 0.488295669245107, 0.472919459698832, 0.564424611143534, 0.405066112938175, 0.580181999133729, 0.497865923446029, 0.492340832264751, 0.300548564690559, 0.574800835942304, 0.671196571821645, 0.458286733858569),
    @("NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries", 0.473929109073491, 0.522198731501057, 0.499553124399261, 0.527639902311216, 0.601436891606279, 0.445671106142509, 0.545641156757768, 0.506066076239562, 0.517911936527468, 0.289612635582178, 0.585085449440321, 0.641734483529905, 0.431123157706663),
    @("International levy on shipping carbon emissions,`nreturned to countries based on population", 0.465444094703709, 0.530655391120507, 0.576537585606553, 0.48644074145873, 0.597737364145906, 0.445997934094592, 0.559112834989109, 0.519427382712729, 0.531035781089868, 0.266453081987289, 0.457115830378939, 0.598511005468507, 0.457869676479742),
    @("Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use", 0.463821638626727, 0.538759689922481, 0.52977181208985, 0.525304380591818, 0.617252086958905, 0.490314971827439, 0.554091313479052, 0.529583432342827, 0.504365824054247, 0.301602146148853, 0.356274513713276, 0.60469307603435, 0.437206606903224),
    @("International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population", 0.384152182838374, 0.42600422832981, 0.454491185458705, 0.413684302919122, 0.443948328445448, 0.39426926452943, 0.4382918870814, 0.42536962433543, 0.408990355395134, 0.242119668148945, 0.339240679877308, 0.54630217154287, 0.368011157501557)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        if ($r -eq 0 -and $c -eq 0) {
            # A1 is already an empty string in the source workbook and is
            # unchanged by this edit - leave it untouched rather than writing
            # an empty value (which would clear the cell instead of keeping
            # it as an explicit empty string).
            continue
        }
        $ws.Cells.Item($r + 1, $c + 1).Value2 = $rowVals[$c]
    }
}

Write-Output "done"
